# Fixed Sprint Backlog Formatting.
# - Remove the stray "D1-Sprint 2" sheet (duplicate/unused Deliverable-1 tab)
# - Rename the remaining Deliverable-2 sprint sheets so numbering is contiguous
#   (D2-Sprint 3 -> D2-Sprint 2, D2-Sprint 4 -> D2-Sprint 3)
# - Reset the stray full-column selection on "D1-Sprint 1" back to a normal cell
# - Leave the new last sheet ("D2-Sprint 3") as the active tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Fix up the selection on the first sheet (was a full-column selection A1:XFD1048576).
$ws1 = $wb.Worksheets.Item("D1-Sprint 1")
$ws1.Activate() | Out-Null
$ws1.Range("A25").Select() | Out-Null

# Delete the extra "D1-Sprint 2" sheet entirely.
$ws2 = $wb.Worksheets.Item("D1-Sprint 2")
$ws2.Delete() | Out-Null

# Renumber the Deliverable-2 sprint sheets.
$ws3 = $wb.Worksheets.Item("D2-Sprint 3")
$ws3.Name = "D2-Sprint 2"

$ws4 = $wb.Worksheets.Item("D2-Sprint 4")
$ws4.Name = "D2-Sprint 3"

# Leave the last sheet selected/active, matching the saved workbook state.
$ws4.Activate() | Out-Null
